# "Version 2." -> "Version 1." (re-reverting the wireframes version bump),
# while preserving the _GoBack bookmark that sits between the " 2"/"." runs
# and the rest of the paragraph's run layout as closely as possible.
#
# Original run layout:
#   r1: "Versi"   r2: "on"   r3: " 2"   <bookmarkStart/End>   r4: "."
# Target run layout:
#   r1: "Version"            r2: " 1."  <bookmarkStart/End>
#
$d = $word.ActiveDocument

# Step 1: "Versi" + "on" (two separate runs) collapse into a single run
# reading "Version". Clear the first run's text, then replace "on" with the
# full word - this leaves exactly one run behind.
$rVersi = $d.Content
$rVersi.Find.Execute("Versi", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$rOn = $d.Content
$rOn.Find.Execute("on", $false, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# Step 2: remove the trailing "." run. Do this *before* touching " 2" so the
# edit never has to span the bookmark that lives between the " 2" run and
# the "." run (an edit spanning the bookmark would delete it).
$periodRange = $d.Content
$periodRange.Find.Execute(".")
if ($periodRange.Find.Found) {
    $periodRange.Delete()
}

# Step 3: " 2" -> " 1." - this run is still entirely before the bookmark, so
# the bookmark is left untouched in between this run and the paragraph end.
$rTwo = $d.Content
$rTwo.Find.Execute(" 2", $false, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)
